# Remove the three "razem" (total) rows from the "nauczyciele" sheet —
# one per year (2021, 2020, 2019). Deleting from the bottom up keeps the
# row numbers of the not-yet-deleted rows stable.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nauczyciele")

$ws.Rows("19:19").Delete()
$ws.Rows("13:13").Delete()
$ws.Rows("7:7").Delete()

# Leave the selection on the (now empty) row right after the data, mirroring
# where Excel parks the cursor after a row deletion of this kind.
$ws.Rows("17:17").Select()
